$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 24: "Move all zeroes to end of array" -----------------------
$ws.Range("A24").Value = "GFG"
$ws.Range("B24").Value = "Move all zeroes to end of array"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "Java"
$ws.Range("D24").Value = 45000
$ws.Range("D24").NumberFormat = $ws.Range("D18").NumberFormat

# ---- Row 25: "Searching an element in a sorted array" ----------------
$ws.Range("A25").Value = "GFG"
$ws.Range("B25").Value = "Searching an element in a sorted array"
$ws.Range("B25").Font.Bold = $true
$ws.Range("B25").Font.Size = 13.5
$ws.Range("B25").WrapText = $false
$ws.Range("B25").VerticalAlignment = -4108
$ws.Range("C25").Value = "Java"
$ws.Range("D25").Value = 45000
$ws.Range("D25").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Rows.Item(25).RowHeight = 18

$ws.Range("C25:D25").Select()
